$wb = $excel.ActiveWorkbook

# --- references sheet: add two new citation/footnote rows ---
$wsRef = $wb.Worksheets.Item("references")
# Enter A12 first so its text becomes shared-string index 20, then A11 so its
# text becomes shared-string index 21 (matches the order newly-added strings
# appear in xl/sharedStrings.xml).
$wsRef.Range("A12").Value = "Klein RJ, Schoenborn CA. Age adjustment using the 2000 projected U.S. population. Healthy People Statistical Notes, no. 20. Hyattsville, Maryland: National Center for Health Statistics. January 2001."
$wsRef.Range("A11").Value = "Table 2, Distribution #1, with under 1 year and 1-4 years combined"
# Reset the selection back to the top-left cell.
$wsRef.Range("A1").Select() | Out-Null

# --- data sheet: move the selection to C14 ---
$wsData = $wb.Worksheets.Item("data")
$wsData.PageSetup.Orientation = 1
$wsData.Range("C14").Select() | Out-Null

# --- drinkAge becomes the active sheet/tab ---
$wsDrink = $wb.Worksheets.Item("drinkAge")
$wsDrink.Activate() | Out-Null
